{"js": "// Office.js (Word JavaScript API) script.\n// Appends the eight new \"final report\" paragraphs (library-system writeup\n// covering login, main menu, borrow, return, reserve, cancel-reserve and\n// renew) after the existing last paragraph (\"This returns the book.\"),\n// and adds a trailing empty paragraph carrying the \"_GoBack\" bookmark that\n// Word leaves at the last edit position.\n\nconst newParagraphTexts = [\n  \"In order to login, the system asks for a username and password. The password is masked to protect the user\\u2019s credentials. The system then takes the username and checks the library\\u2019s vectors for students, teachers, and librarians. The system also keeps track of the type of user as well as the index of the user in that specific user vector. This allows the system to use functions on the correct user. After the system has logged the user in, it displays the menu and displays the options for each user. For the inputs on the main menu for each user, it uses an overloaded function in the user class. That means that there is a base, generic menu for a user and a different menu for the derived classes. \",\n  \"In order for the user to change do anything in the system, they input a number signifying what action they want the system to do. The system also checks what type of user they are in order to accept different types of commands for a reader or a librarian. Each command also has a different case for if the user is a student or a teacher.\",\n  \"Borrowing a book requires the user to input a value for the ISBN of the book they wish to borrow. The system checks if the ISBN is of a book that exists and ends the function if the ISBN is invalid. It also checks through all the copies for if it has already been borrowed by the user. It also checks if the user has any overdue books to stop them from borrowing. It then checks for any available copy to borrow and sets the appropriate information in the copy of the book as well as adds it into the user\\u2019s borrowed vector. It also works if the user is the first one on the reserve list.\",\n  \"Returning a book requires the user to input a value for the id of the copy they wish to return. It checks if book being returned is overdue and gives a penalty to the user. If the penalties are a multiple of 5, it reduces the maxCopies by 1. It resets all the information in the copy and the user. After all that, it gives a poll to the user on whether they enjoyed the book or not.\",\n  \"Reserving a book requires the user to input a value for the ISBN of the book they wish to reserve. It checks if the ISBN for validity and also checks if a user has reserved the book to which it stops the command. It adds the user to the reserve list and puts the user in the copy if they are the only reserve there. \",\n  \"Cancelling a reserve requires the user to input a value for the ISBN of the book they wish to cancel the reserve. It checks for the validity of the ISBN and erases the user from the reserve list. If there is no reserve left, it returns the reserve in the copy to the default state.\",\n  \"Renewing requires the user to input an id for the book they wish to renew. The system checks for if it is inside their borrowed copies vector and adds another instance of the maxBorrowingPeriod. \"\n];\n\n// Find the last paragraph in the body (\"This returns the book.\") to anchor\n// the insertion point.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert each new paragraph after the current anchor, matching the\n// first-line-indent formatting used by the rest of the body text.\nfor (const text of newParagraphTexts) {\n  const inserted = anchor.insertParagraph(text, Word.InsertLocation.after);\n  inserted.firstLineIndent = 36; // 720 twips == 0.5in == 36pt\n  anchor = inserted;\n  await context.sync();\n}\n\n// Final trailing empty paragraph (same indent) holding the \"_GoBack\"\n// bookmark Word drops at the last edited location before saving.\nconst trailing = anchor.insertParagraph(\"\", Word.InsertLocation.after);\ntrailing.firstLineIndent = 36;\nawait context.sync();\n\nconst trailingRange = trailing.getRange();\ntrailingRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Appends the eight new \"final report\" paragraphs (library-system writeup\n# covering login, main menu, borrow, return, reserve, cancel-reserve and\n# renew) after the existing last paragraph (\"This returns the book.\"),\n# and adds a trailing empty paragraph carrying the \"_GoBack\" bookmark that\n# Word leaves at the last edit position.\n\n$d = $word.ActiveDocument\n\n$newParagraphTexts = @(\n\"In order to login, the system asks for a username and password. The password is masked to protect the user\u2019s credentials. The system then takes the username and checks the library\u2019s vectors for students, teachers, and librarians. The system also keeps track of the type of user as well as the index of the user in that specific user vector. This allows the system to use functions on the correct user. After the system has logged the user in, it displays the menu and displays the options for each user. For the inputs on the main menu for each user, it uses an overloaded function in the user class. That means that there is a base, generic menu for a user and a different menu for the derived classes. \",\n\"In order for the user to change do anything in the system, they input a number signifying what action they want the system to do. The system also checks what type of user they are in order to accept different types of commands for a reader or a librarian. Each command also has a different case for if the user is a student or a teacher.\",\n\"Borrowing a book requires the user to input a value for the ISBN of the book they wish to borrow. The system checks if the ISBN is of a book that exists and ends the function if the ISBN is invalid. It also checks through all the copies for if it has already been borrowed by the user. It also checks if the user has any overdue books to stop them from borrowing. It then checks for any available copy to borrow and sets the appropriate information in the copy of the book as well as adds it into the user\u2019s borrowed vector. It also works if the user is the first one on the reserve list.\",\n\"Returning a book requires the user to input a value for the id of the copy they wish to return. It checks if book being returned is overdue and gives a penalty to the user. If the penalties are a multiple of 5, it reduces the maxCopies by 1. It resets all the information in the copy and the user. After all that, it gives a poll to the user on whether they enjoyed the book or not.\",\n\"Reserving a book requires the user to input a value for the ISBN of the book they wish to reserve. It checks if the ISBN for validity and also checks if a user has reserved the book to which it stops the command. It adds the user to the reserve list and puts the user in the copy if they are the only reserve there. \",\n\"Cancelling a reserve requires the user to input a value for the ISBN of the book they wish to cancel the reserve. It checks for the validity of the ISBN and erases the user from the reserve list. If there is no reserve left, it returns the reserve in the copy to the default state.\",\n\"Renewing requires the user to input an id for the book they wish to renew. The system checks for if it is inside their borrowed copies vector and adds another instance of the maxBorrowingPeriod. \"\n)\n\nforeach ($t in $newParagraphTexts) {\n    # Re-query the last paragraph each time so the inserted paragraph mark\n    # and new text land in the right place.\n    $lastPara = $d.Paragraphs.Last\n    $r = $lastPara.Range\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n\n    $newPara = $d.Paragraphs.Last\n    $r2 = $newPara.Range\n    $r2.Collapse(0)\n    $r2.InsertAfter($t)\n}\n\n# Trailing empty paragraph (inherits the same first-line indent) holding\n# the \"_GoBack\" bookmark Word drops at the last edited location before\n# saving.\n$lastPara = $d.Paragraphs.Last\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$bmPara = $d.Paragraphs.Last\n$bmRange = $bmPara.Range\n$bmRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n"}
